# ---------------------------------------------------------------------------
# Glenn Dalbey resume edit:
#   - "74th open / 94th closed of 1,134 teams" -> "Top 8% of 1,134 teams"
#     (summary paragraph + project subtitle)
#   - Missing Persons bullet: 44.75 sigma -> "up to 46.86" sigma
#   - New "OE-OS (In Progress)" project block inserted before "AI Homelab"
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$bullet = [char]0x2022

function Replace-ParagraphXml($matchText, $innerXml) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "$matchText*") {
            $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $p.Range.InsertXML($pkg)
            return $true
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1) PROFESSIONAL SUMMARY paragraph: update the Kaggle placement parenthetical
# ---------------------------------------------------------------------------
$summaryXml = '<w:p><w:pPr><w:spacing w:after="120"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Data Science professional with MS in Data Science and proven expertise in competitive machine learning and deep learning systems. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Kaggle Bronze Medalist</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> in NFL Big Data Bowl 2026 (Top 8% of 1,134 teams) with 847+ experiments across 15+ neural network architectures. Trained 105 3D medical imaging models and deployed production healthcare AI achieving 93.8% accuracy. Expert in spatial-temporal modeling, trajectory prediction, ensemble methods, and multi-modal AI. Strong foundation in systematic ML experimentation, advanced feature engineering, and production deployment.</w:t></w:r>' +
    '</w:p>'
Replace-ParagraphXml "Data Science professional" $summaryXml | Out-Null

# ---------------------------------------------------------------------------
# 2) "Player Trajectory Prediction | ..." project subtitle
# ---------------------------------------------------------------------------
$subtitleXml = '<w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Player Trajectory Prediction | Top 8% of 1,134 teams</w:t></w:r></w:p>'
Replace-ParagraphXml "Player Trajectory Prediction" $subtitleXml | Out-Null

# ---------------------------------------------------------------------------
# 3) Missing Persons bullet: sigma significance figure
# ---------------------------------------------------------------------------
$sigmaXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $bullet + ' Analyzed 41,200 cases across 101 years identifying trafficking corridors at up to 46.86' + [char]0x3C3 + ' significance</w:t></w:r></w:p>'
Replace-ParagraphXml "$bullet Analyzed 41,200 cases" $sigmaXml | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert the new "OE-OS (In Progress)" project block right after the
#    "Built 7-page interactive Streamlit dashboard..." bullet and before the
#    "AI Homelab & Active Memory Network" heading.
#
#    Strategy: drop placeholder paragraphs via InsertAfter (keeps paragraph
#    count/order correct without inheriting neighboring run formatting),
#    then stamp each placeholder with its exact target OOXML via
#    Range.InsertXML so formatting matches the target precisely.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$bullet Built 7-page interactive Streamlit dashboard*") {
        $anchor = $p.Range.Duplicate
        $anchor.Collapse(0)
        $anchor.InsertAfter("OEOSPLACEHOLDER1`rOEOSPLACEHOLDER2`rOEOSPLACEHOLDER3`rOEOSPLACEHOLDER4`rOEOSPLACEHOLDER5`rOEOSPLACEHOLDER6`r")
        break
    }
}

Replace-ParagraphXml "OEOSPLACEHOLDER1" '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">OE-OS (In Progress)</w:t></w:r></w:p>' | Out-Null

Replace-ParagraphXml "OEOSPLACEHOLDER2" '<w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Distributed AI Orchestration Platform | Python / FastAPI</w:t></w:r></w:p>' | Out-Null

Replace-ParagraphXml "OEOSPLACEHOLDER3" ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $bullet + ' Three-tier LLM routing (local Ollama to cheap API to premium models) routing ~80% of requests to free local models</w:t></w:r></w:p>') | Out-Null

Replace-ParagraphXml "OEOSPLACEHOLDER4" ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $bullet + ' Triple-layer RAG memory: BM25 over 5M+ chunks, ChromaDB semantic search, Redis session cache</w:t></w:r></w:p>') | Out-Null

Replace-ParagraphXml "OEOSPLACEHOLDER5" ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $bullet + ' 18 MCP-compatible tools and multi-agent sandbox with 4 LLM personas at zero API cost</w:t></w:r></w:p>') | Out-Null

Replace-ParagraphXml "OEOSPLACEHOLDER6" ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="100"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $bullet + ' 4,200+ lines of async Python on FastAPI for a private multi-node GPU cluster</w:t></w:r></w:p>') | Out-Null

# The InsertAfter splice above can strip the (unneeded-but-present-in-source)
# xml:space="preserve" from the untouched "AI Homelab..." run that follows
# it; restore that paragraph's exact original markup so it stays byte-for-
# byte identical to the source.
Replace-ParagraphXml "AI Homelab" '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">AI Homelab &amp; Active Memory Network</w:t></w:r></w:p>' | Out-Null

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
